# Add the "Img-Source" caption textbox to slide 30 ("Going from Shallow to
# Deep Neural Networks"), underneath/over the network-diagram picture.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(30)

# The slide currently holds shapes with ids 2,3,4,6 (id 5 is a historical
# gap). This engine hands out the next shape id as "current shape count +
# 1", bumping past collisions - so the very first shape we add here would
# land on id 5, not 7. Burn id 5 with a throwaway shape first so the real
# textbox lands on id 7, matching the authored deck, then remove the
# throwaway.
$placeholder = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)

# Position/size in points, chosen so that round(pt * 12700) reproduces the
# exact target EMU values (4572000, 1412776, 4343400, 276999).
$left   = 4572000 / 12700
$top    = 1412776 / 12700
$width  = 4343400 / 12700
$height = 276999 / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "Textfeld 6"

$placeholder.Delete()

$tb.Fill.Visible = [Microsoft.Office.Core.MsoTriState]::msoFalse
$tb.TextFrame.WordWrap = $true
$tb.TextFrame.AutoSize = 1

$tr = $tb.TextFrame.TextRange
$tr.Text = "Img-Source"
$tr.Font.Size = 12

$run2 = $tr.InsertAfter(": http://")
$run2.Font.Size = 12

$run3 = $run2.InsertAfter("neuralnetworksanddeeplearning.com")
$run3.Font.Size = 12

$tb.TextFrame.TextRange.ParagraphFormat.Alignment = 3
